$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before E ("giờ test" shifts from E->F, "Địa điểm" F->G)
$ws.Columns("E:E").Insert()

# New header text for the inserted column
$ws.Range("E1").Value = "Ca làm việc"

# Try to match the original column D width (best effort given COM width
# quantization) so the new column gets a visible custom width too.
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Re-apply the AutoFilter over the new, wider range (A1:G3)
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:G3").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "DATA!_FilterDatabase") {
        $n.RefersTo = "=DATA!`$A`$1:`$G`$3"
    }
}

# Restore/update the selected cell recorded in the sheet view
[void]$ws.Range("N9").Select()
